# Generate Report for handoff
#
# - rename the in-flight handoff file from the old UUID to a new UUID
#   (c36e646e-08fb-4c64-a4bb-c4437a834597) with a new content hash
#   (67715268c462517be7404f981a7134af577bbc21), refreshing the handoff
#   timestamps
# - record that the PREVIOUS handoff (the old UUID's .md) failed its
#   handoff transform, by adding a new "Handoff transform failed" row
#   for it (previously "Ready for handoff")
# - add a new row for ".localization-config" (status "Not to be
#   localized"), which is what used to sit in the row that is now
#   used by the failed handoff.

$wb = $excel.ActiveWorkbook

$oldUuid  = "d4c99864-2751-42f0-8ed8-8cabfa458d72"
$newUuid  = "c36e646e-08fb-4c64-a4bb-c4437a834597"
$failUuid = "b74e7a9b-d654-4e39-be77-9c6af37f6cce"
$newHash  = "67715268c462517be7404f981a7134af577bbc21"

$newMdName   = "$newUuid.md"
$failMdName  = "$failUuid.md"
$cfgName     = ".localization-config"

$zhXlfName = "$newUuid.$newHash.zh-cn.xlf"
$deXlfName = "$newUuid.$newHash.de-de.xlf"

$zhDatetime = "2016-02-15 04:05:39"
$deDatetime = "2016-02-15 04:05:52"

$baseUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/38781ae7fce093b92ba897f08344ebe45d2eecd7"
$zhTargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/122d9d1165d67da455bdc4aefa27096d7ff8aa9a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlfName"
$deTargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1de59aaeee4efa3a287ab9b172040a914ff9caa3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlfName"

$newMdUrl   = "$baseUrl/e2e/$newMdName"
$failMdUrl  = "$baseUrl/e2e/$failMdName"
$cfgUrl     = "$baseUrl/$cfgName"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 3 used to be the ".localization-config" row; it now becomes the
# "failed handoff" row for the old UUID, and a brand-new row 4 takes
# over the ".localization-config" entry.
$ws.Range("B3").Value = "Handoff transform failed"
$ws.Range("C3").Value = "Handoff transform failed"

$ws.Range("B4").Value = "Not to be localized"
$ws.Range("C4").Value = "Not to be localized"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", $newMdName)
$ws.Hyperlinks.Add($ws.Range("A3"), $failMdUrl, "", "", $failMdName)
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", $cfgName)

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("D2").Value = $zhDatetime

$ws.Range("B3").Value = "Handoff transform failed"

$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", $newMdName)
$ws.Hyperlinks.Add($ws.Range("C2"), $zhTargetUrl, "", "", $zhXlfName)
$ws.Hyperlinks.Add($ws.Range("A3"), $failMdUrl, "", "", $failMdName)
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", $cfgName)

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("D2").Value = $deDatetime

$ws.Range("B3").Value = "Handoff transform failed"

$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, "", "", $newMdName)
$ws.Hyperlinks.Add($ws.Range("C2"), $deTargetUrl, "", "", $deXlfName)
$ws.Hyperlinks.Add($ws.Range("A3"), $failMdUrl, "", "", $failMdName)
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", $cfgName)
